{"js": "// Resume edit: add a new \"Experienced in Lua game scripting\" bullet at the\n// top of the \"TECHNICAL SKILLS AND QUALIFICATIONS\" list, and move the\n// \"_GoBack\" bookmark (which Word drops at the location of the most recent\n// edit) from its old spot in \"Individual Python project\" onto the newly\n// inserted bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the first bullet of the technical-skills list; this is the\n// paragraph whose text starts the shifted block in the diff.\nconst targetText =\n  \"Experienced in custom back-end development across platforms (mobile & desktop & web)\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Experienced in custom back-end...' paragraph\");\n}\n\n// Insert the new bullet immediately above it; it inherits the list\n// paragraph's style/format (bullet numbering, indent, fonts, etc.).\ntarget.insertParagraph(\"Experienced in Lua game scripting\", \"Before\");\nawait context.sync();\n\n// Re-fetch paragraphs so the newly minted paragraph's range collapses\n// correctly (a freshly-returned insert proxy can carry a stale anchor).\nconst paragraphsAfterInsert = body.paragraphs;\nparagraphsAfterInsert.load(\"items/text\");\nawait context.sync();\n\nlet newPara = null;\nfor (let i = 0; i < paragraphsAfterInsert.items.length; i++) {\n  if (paragraphsAfterInsert.items[i].text === \"Experienced in Lua game scripting\") {\n    newPara = paragraphsAfterInsert.items[i];\n    break;\n  }\n}\nif (!newPara) {\n  throw new Error(\"Could not find the newly inserted 'Experienced in Lua game scripting' paragraph\");\n}\n\n// Move the \"_GoBack\" bookmark onto the end of the freshly inserted bullet.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst newRange = newPara.getRange(\"End\");\nnewRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$target = $null\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Experienced in custom back-end development across platforms (mobile & desktop & web)\") {\n        $target = $p\n        $targetIndex = $i\n        break\n    }\n}\nWrite-Output \"targetIndex=$targetIndex\"\n\n$target.Range.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs.Item($targetIndex)\nWrite-Output \"newPara text=[$($newPara.Range.Text)]\"\n$newPara.Range.Text = \"Experienced in Lua game scripting\"\nWrite-Output \"after set: [$($newPara.Range.Text)]\"\n"}
